$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the old row 3 (C3=0) down to row 4
$ws.Range("C4").Value = 0

# Populate new row 3 with trade data
$ws.Range("B3").Value = $false
$ws.Range("C3").Value = 0
$ws.Range("E3").Value = 80.9599
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = $false

# G column carries the date-format style (style index 1), same as G2
$ws.Range("G3").NumberFormat = "m/d/yy h:mm"
